$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows, restricted to the used A:B columns, pushing the
# existing data down (rows shift 2 -> 4, 3 -> 5, etc.).
$ws.Range("A2:B3").Insert(-4121)

# Copy formatting from the row that now holds the original "row 2" data
# (originally row 2, now row 4) into the two newly inserted rows.
$ws.Range("A4:B4").Copy()
$ws.Range("A2:B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the new rows' values.
$ws.Range("A2").Value = 45750
$ws.Range("B2").Value = "Annegudu"

$ws.Range("A3").Value = 45751
$ws.Range("B3").Value = "Annegudu"

# Match the row height of the template row (the thick-bottom border under
# each data row makes Excel bump the row height slightly).
$ws.Rows.Item(2).RowHeight = $ws.Rows.Item(4).RowHeight
$ws.Rows.Item(3).RowHeight = $ws.Rows.Item(4).RowHeight

# Leave the selection on A4, matching the saved selection state.
$ws.Range("A4").Select()
